# The author added a new weekly price observation for Berenjena at
# Vega Monumental Concepción. This inserts a new data row right after the
# existing row 107 (i.e. at row 108), pushing the previous rows 108-144
# down to 109-145, and fills the new row with the reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108; everything below shifts down by one.
$ws.Rows(108).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A108").Value = 11
$ws.Range("B108").Value = "Vega Monumental Concepción"
$ws.Range("C108").Value = "Bíobío"
$ws.Range("D108").Value = 45009
$ws.Range("E108").Value = 8
$ws.Range("F108").Value = 100112001
$ws.Range("G108").Value = "Berenjena"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 100
$ws.Range("K108").Value = 8000
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 8500
$ws.Range("N108").Value = "$/caja 60 unidades"
$ws.Range("O108").Value = "Región de Arica y Parinacota"
$ws.Range("P108").Value = 142
$ws.Range("Q108").Value = 60
$ws.Range("R108").Value = "Hortaliza"
